# Update for release to deply 0.1.1
$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsInclude = $wb.Worksheets.Item("Include from NMDP Transplant ")

# Bump the version number shown on the Metadata sheet
$wsMeta.Range("B3").Value = "0.1.1"

# Refresh the publication date
$wsMeta.Range("B8").Value = "2024-11-11T17:53:38-06:00"

# Insert a new "Jurisdiction" property row before "Description"
$wsMeta.Rows.Item(11).Insert()

# Match formatting of the surrounding data rows
$wsMeta.Range("A10:B10").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$wsMeta.Range("A11").Value = "Jurisdiction"
$wsMeta.Range("B11").Value = ""

# Rename the include tab to match the new release
$wsInclude.Name = "Include #0"

$wb.Save()
